$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_1_9_24"
$ws.Range("B2").Value = 0.9918592411180278
$ws.Range("C2").Value = 0.9502644847130644
$ws.Range("D2").Value = 0.9695432369876992
$ws.Range("E2").Value = 0.9601327186165663
$ws.Range("F2").Value = 0.9817779449189653
$ws.Range("G2").Value = 3.748475519549361
$ws.Range("H2").Value = 2.649151507048876
$ws.Range("I2").Value = 3.231148589633617

$ws.Range("A3").Value = "model_1_9_23"
$ws.Range("B3").Value = 0.9918845064403244
$ws.Range("C3").Value = 0.950396423850235
$ws.Range("D3").Value = 0.9696767725688438
$ws.Range("E3").Value = 0.9602650877662967
$ws.Range("F3").Value = 0.9787309395277349
$ws.Range("G3").Value = 3.738531506244158
$ws.Range("H3").Value = 2.63753648460243
$ws.Range("I3").Value = 3.220420384031915

$ws.Range("A4").Value = "model_1_9_22"
$ws.Range("B4").Value = 0.9919121459727243
$ws.Range("C4").Value = 0.9505416529557069
$ws.Range("D4").Value = 0.96982425594215
$ws.Range("E4").Value = 0.9604110805977357
$ws.Range("F4").Value = 0.975397603691181
$ws.Range("G4").Value = 3.72758585214875
$ws.Range("H4").Value = 2.624708272999606
$ws.Range("I4").Value = 3.20858801134355

$ws.Range("A5").Value = "model_1_9_21"
$ws.Range("B5").Value = 0.9919423423573986
$ws.Range("C5").Value = 0.9507013895480402
$ws.Range("D5").Value = 0.9699872202100419
$ws.Range("E5").Value = 0.96057208102787
$ws.Range("F5").Value = 0.9717559107090704
$ws.Range("G5").Value = 3.715546795099016
$ws.Range("H5").Value = 2.61053352187105
$ws.Range("I5").Value = 3.195539308379446

$ws.Range("A6").Value = "model_1_9_20"
$ws.Range("B6").Value = 0.9919752642863411
$ws.Range("C6").Value = 0.9508770733925053
$ws.Range("D6").Value = 0.9701672909510058
$ws.Range("E6").Value = 0.9607494889003967
$ws.Range("F6").Value = 0.9677855162767417
$ws.Range("G6").Value = 3.702305822599621
$ws.Range("H6").Value = 2.594870837211901
$ws.Range("I6").Value = 3.181160821128432

$ws.Range("A7").Value = "model_1_9_19"
$ws.Range("B7").Value = 0.9920110508198798
$ws.Range("C7").Value = 0.9510697853031066
$ws.Range("D7").Value = 0.9703659353627258
$ws.Range("E7").Value = 0.960944690600783
$ws.Range("F7").Value = 0.9634696496772268
$ws.Range("G7").Value = 3.687781475660669
$ws.Range("H7").Value = 2.57759260109528
$ws.Range("I7").Value = 3.16534018633692

$ws.Range("A8").Value = "model_1_9_18"
$ws.Range("B8").Value = 0.9920498738239156
$ws.Range("C8").Value = 0.9512809267892706
$ws.Range("D8").Value = 0.9705852179294737
$ws.Range("E8").Value = 0.9611594248932196
$ws.Range("F8").Value = 0.9587875838317037
$ws.Range("G8").Value = 3.671868125060388
$ws.Range("H8").Value = 2.558519243170295
$ws.Range("I8").Value = 3.147936481291696

$ws.Range("A9").Value = "model_1_9_17"
$ws.Range("B9").Value = 0.9920917817834813
$ws.Range("C9").Value = 0.9515119535864197
$ws.Range("D9").Value = 0.9708264377034345
$ws.Range("E9").Value = 0.9613950030471756
$ws.Range("F9").Value = 0.9537334714308026
$ws.Range("G9").Value = 3.65445605466205
$ws.Range("H9").Value = 2.537537770928481
$ws.Range("I9").Value = 3.128843430712633

$ws.Range("A10").Value = "model_1_9_16"
$ws.Range("B10").Value = 0.9921368804112949
$ws.Range("C10").Value = 0.9517640961265175
$ws.Range("D10").Value = 0.9710921740342916
$ws.Range("E10").Value = 0.961653301989789
$ws.Range("F10").Value = 0.9482945634892516
$ws.Range("G10").Value = 3.635452529041763
$ws.Range("H10").Value = 2.51442383064916
$ws.Range("I10").Value = 3.107908914107347

$ws.Range("A11").Value = "model_1_9_15"
$ws.Range("B11").Value = 0.9921851602723137
$ws.Range("C11").Value = 0.9520385902972514
$ws.Range("D11").Value = 0.9713841508858527
$ws.Range("E11").Value = 0.9619359316112902
$ws.Range("F11").Value = 0.942471997875977
$ws.Range("G11").Value = 3.614764401587593
$ws.Range("H11").Value = 2.489027470700335
$ws.Range("I11").Value = 3.085002453691368

$ws.Range("A12").Value = "model_1_9_14"
$ws.Range("B12").Value = 0.9922365559796579
$ws.Range("C12").Value = 0.9523368301823291
$ws.Range("D12").Value = 0.9717046313850156
$ws.Range("E12").Value = 0.9622445976653697
$ws.Range("F12").Value = 0.9362736602681958
$ws.Range("G12").Value = 3.592286602740686
$ws.Range("H12").Value = 2.461151842650357
$ws.Range("I12").Value = 3.059985802174193

$ws.Range("A13").Value = "model_1_9_13"
$ws.Range("B13").Value = 0.9922908665882192
$ws.Range("C13").Value = 0.9526595953006751
$ws.Range("D13").Value = 0.972055905319455
$ws.Range("E13").Value = 0.9625808718006241
$ws.Range("F13").Value = 0.929723784705782
$ws.Range("G13").Value = 3.567960381574492
$ws.Range("H13").Value = 2.430597779093724
$ws.Range("I13").Value = 3.032731581165056

$ws.Range("A14").Value = "model_1_9_12"
$ws.Range("B14").Value = 0.9923478042025029
$ws.Range("C14").Value = 0.9530077308058192
$ws.Range("D14").Value = 0.9724403433307717
$ws.Range("E14").Value = 0.9629464607445434
$ws.Range("F14").Value = 0.9228570914710902
$ws.Range("G14").Value = 3.541722040401391
$ws.Range("H14").Value = 2.397159079891344
$ws.Range("I14").Value = 3.003101464449319

$ws.Range("A15").Value = "model_1_9_11"
$ws.Range("B15").Value = 0.9924068012480828
$ws.Range("C15").Value = 0.9533816858640696
$ws.Range("D15").Value = 0.9728599865715801
$ws.Range("E15").Value = 0.9633424246684117
$ws.Range("F15").Value = 0.9157420302089264
$ws.Range("G15").Value = 3.51353772637195
$ws.Range("H15").Value = 2.360658204096974
$ws.Range("I15").Value = 2.971009527659148

$ws.Range("A16").Value = "model_1_9_10"
$ws.Range("B16").Value = 0.9924671681729413
$ws.Range("C16").Value = 0.9537813086602231
$ws.Range("D16").Value = 0.9733172307244848
$ws.Range("E16").Value = 0.9637701435017164
$ws.Range("F16").Value = 0.9084617610979073
$ws.Range("G16").Value = 3.483418881522486
$ws.Range("H16").Value = 2.32088677348672
$ws.Range("I16").Value = 2.936343930782848

$ws.Range("A17").Value = "model_1_9_9"
$ws.Range("B17").Value = 0.9925277806582722
$ws.Range("C17").Value = 0.9542057313018585
$ws.Range("D17").Value = 0.9738139581072293
$ws.Range("E17").Value = 0.964229974410311
$ws.Range("F17").Value = 0.9011518773208204
$ws.Range("G17").Value = 3.451430917329614
$ws.Range("H17").Value = 2.277681062687502
$ws.Range("I17").Value = 2.899075726375193

$ws.Range("A18").Value = "model_1_9_0"
$ws.Range("B18").Value = 0.9925851208007793
$ws.Range("C18").Value = 0.9582011591958821
$ws.Range("D18").Value = 0.9800732475449591
$ws.Range("E18").Value = 0.96935812017191
$ws.Range("F18").Value = 0.8942366390625439
$ws.Range("G18").Value = 3.150302768471241
$ws.Range("H18").Value = 1.733243492604293
$ws.Range("I18").Value = 2.483451676526844

$ws.Range("A19").Value = "model_1_9_8"
$ws.Range("B19").Value = 0.9925871439413033
$ws.Range("C19").Value = 0.9546530723527547
$ws.Range("D19").Value = 0.9743518742622788
$ws.Range("E19").Value = 0.9647218305471579
$ws.Range("F19").Value = 0.8939926477129989
$ws.Range("G19").Value = 3.41771563422659
$ws.Range("H19").Value = 2.230892722369129
$ws.Range("I19").Value = 2.859211953182577

$ws.Range("A20").Value = "model_1_9_7"
$ws.Range("B20").Value = 0.9926432966750481
$ws.Range("C20").Value = 0.9551203454900471
$ws.Range("D20").Value = 0.9749324020922026
$ws.Range("E20").Value = 0.9652450613812377
$ws.Range("F20").Value = 0.8872206113049849
$ws.Range("G20").Value = 3.382498105947678
$ws.Range("H20").Value = 2.180397987426184
$ws.Range("I20").Value = 2.816805335201038

$ws.Range("A21").Value = "model_1_9_1"
$ws.Range("B21").Value = 0.99269132480108
$ws.Range("C21").Value = 0.957894383361084
$ws.Range("D21").Value = 0.9792580680095936
$ws.Range("E21").Value = 0.968795412729386
$ws.Range("F21").Value = 0.8814284050060885
$ws.Range("G21").Value = 3.173423906356226
$ws.Range("H21").Value = 1.804148404388801
$ws.Range("I21").Value = 2.529057779982981

$ws.Range("A22").Value = "model_1_9_6"
$ws.Range("B22").Value = 0.9926935663507477
$ws.Range("C22").Value = 0.9556028603399629
$ws.Range("D22").Value = 0.9755560885239185
$ws.Range("E22").Value = 0.9657976056514989
$ws.Range("F22").Value = 0.8811580734487909
$ws.Range("G22").Value = 3.346131837451346
$ws.Range("H22").Value = 2.126149285755591
$ws.Range("I22").Value = 2.772022932749422

$ws.Range("A23").Value = "model_1_9_5"
$ws.Range("B23").Value = 0.9927344230833944
$ws.Range("C23").Value = 0.9560937380593768
$ws.Range("D23").Value = 0.9762225115244789
$ws.Range("E23").Value = 0.9663757850295182
$ws.Range("F23").Value = 0.876230739327297
$ws.Range("G23").Value = 3.30913527465915
$ws.Range("H23").Value = 2.068183326091599
$ws.Range("I23").Value = 2.725162865621326

$ws.Range("A24").Value = "model_1_9_2"
$ws.Range("B24").Value = 0.9927483882702133
$ws.Range("C24").Value = 0.9575038838749522
$ws.Range("D24").Value = 0.9784550318947349
$ws.Range("E24").Value = 0.9681975708002003
$ws.Range("F24").Value = 0.8745465336390753
$ws.Range("G24").Value = 3.202855143887731
$ws.Range("H24").Value = 1.87399707258224
$ws.Range("I24").Value = 2.577511450242894

$ws.Range("A25").Value = "model_1_9_4"
$ws.Range("B25").Value = 0.9927613891325676
$ws.Range("C25").Value = 0.9565835967662056
$ws.Range("D25").Value = 0.9769303360675665
$ws.Range("E25").Value = 0.9669744900747903
$ws.Range("F25").Value = 0.8729786257683861
$ws.Range("G25").Value = 3.272215512996028
$ws.Range("H25").Value = 2.006616229999035
$ws.Range("I25").Value = 2.676639241849943

$ws.Range("A26").Value = "model_1_9_3"
$ws.Range("B26").Value = 0.992768544220511
$ws.Range("C26").Value = 0.9570593717087106
$ws.Range("D26").Value = 0.9776762868574322
$ws.Range("E26").Value = 0.9675854699251453
$ws.Range("F26").Value = 0.87211572003214
$ws.Range("G26").Value = 3.236357219088632
$ws.Range("H26").Value = 1.941732885096001
$ws.Range("I26").Value = 2.627120774242835
